$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account-number column (C) holds long digit strings that must stay text,
# not be coerced to scientific-notation numbers.
$ws.Range("C2:C7").NumberFormat = "@"

# Row 2 - NASIRI HASNAA (first entry)
$ws.Range("A2").Value = "NASIRI HASNAA"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "546576878798989898090090"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "CIH"
$ws.Range("F2").Value = "Logement de fonction"
$ws.Range("G2").Value = "905/LF/TADLA OUARDIGHA ZAYANE"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 19999.98
$ws.Range("J2").Value = 1999.98
$ws.Range("K2").Value = 18000

# Row 3 - NASIRI HASNAA (second entry)
$ws.Range("A3").Value = "NASIRI HASNAA"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "546576878798989898090090"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "CIH"
$ws.Range("F3").Value = "Logement de fonction"
$ws.Range("G3").Value = "905/LF/TADLA OUARDIGHA ZAYANE"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 19999.98
$ws.Range("J3").Value = 999.99
$ws.Range("K3").Value = 18999.99

# Row 4 - MOHAMED BADRANE
$ws.Range("A4").Value = "MOHAMED BADRANE"
$ws.Range("B4").Value = "I83603"
$ws.Range("C4").Value = "225400000805987601012173"
$ws.Range("D4").Value = "KHOURIBGA"
$ws.Range("E4").Value = "CA"
$ws.Range("F4").Value = "Point de vente"
$ws.Range("G4").Value = "605/KHOURIBGA NAHDA"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 7500
$ws.Range("J4").Value = 375
$ws.Range("K4").Value = 7125

# Row 5 - ZERNAKH ABDELLAH
$ws.Range("A5").Value = "ZERNAKH ABDELLAH"
$ws.Range("B5").Value = "IB19558"
$ws.Range("C5").Value = "145101211406073828000084"
$ws.Range("D5").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E5").Value = "BP"
$ws.Range("F5").Value = "Point de vente"
$ws.Range("G5").Value = "052/FKIH BEN SALEH"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 11000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 11000

# Row 6 - NOUBAIL MOUNTASSIR
$ws.Range("A6").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B6").Value = "Q251990"
$ws.Range("C6").Value = "007400000313200019604463"
$ws.Range("D6").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E6").Value = "AWB"
$ws.Range("F6").Value = "Direction régionale"
$ws.Range("G6").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H6").Value = "mensuelle"
$ws.Range("I6").Value = 6750
$ws.Range("J6").Value = 675
$ws.Range("K6").Value = 6075

# Row 7 - NOUBAIL MOHAMMED
$ws.Range("A7").Value = "NOUBAIL MOHAMMED"
$ws.Range("B7").Value = "IR801997"
$ws.Range("C7").Value = "007400000313200019604463"
$ws.Range("D7").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E7").Value = "AWB"
$ws.Range("F7").Value = "Direction régionale"
$ws.Range("G7").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H7").Value = "mensuelle"
$ws.Range("I7").Value = 6750
$ws.Range("J7").Value = 675
$ws.Range("K7").Value = 6075

# Row 8 - totals row
$ws.Range("A8").Value = " "
$ws.Range("B8").Value = " "
$ws.Range("C8").Value = " "
$ws.Range("D8").Value = " "
$ws.Range("E8").Value = " "
$ws.Range("F8").Value = " "
$ws.Range("G8").Value = " "
$ws.Range("H8").Value = " "
$ws.Range("I8").Value = 71999.96
$ws.Range("J8").Value = 4724.97
$ws.Range("K8").Value = 67274.99
